$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 2..97 (row 1 is the header row, left unchanged)
# Column A = Timestamp (date serial), Column B = Actual Production (MW)
$rowsData = @(
    [PSCustomObject]@{ Row = 2; A = 45890.01041666666; B = 894 }
    [PSCustomObject]@{ Row = 3; A = 45890.02083333334; B = 838 }
    [PSCustomObject]@{ Row = 4; A = 45890.03125; B = 818 }
    [PSCustomObject]@{ Row = 5; A = 45890.04166666666; B = 797 }
    [PSCustomObject]@{ Row = 6; A = 45890.05208333334; B = 744 }
    [PSCustomObject]@{ Row = 7; A = 45890.0625; B = 692 }
    [PSCustomObject]@{ Row = 8; A = 45890.07291666666; B = 650 }
    [PSCustomObject]@{ Row = 9; A = 45890.08333333334; B = 612 }
    [PSCustomObject]@{ Row = 10; A = 45890.09375; B = 586 }
    [PSCustomObject]@{ Row = 11; A = 45890.10416666666; B = 547 }
    [PSCustomObject]@{ Row = 12; A = 45890.11458333334; B = 500 }
    [PSCustomObject]@{ Row = 13; A = 45890.125; B = 454 }
    [PSCustomObject]@{ Row = 14; A = 45890.13541666666; B = 423 }
    [PSCustomObject]@{ Row = 15; A = 45890.14583333334; B = 391 }
    [PSCustomObject]@{ Row = 16; A = 45890.15625; B = 373 }
    [PSCustomObject]@{ Row = 17; A = 45890.16666666666; B = 373 }
    [PSCustomObject]@{ Row = 18; A = 45890.17708333334; B = 379 }
    [PSCustomObject]@{ Row = 19; A = 45890.1875; B = 396 }
    [PSCustomObject]@{ Row = 20; A = 45890.19791666666; B = 405 }
    [PSCustomObject]@{ Row = 21; A = 45890.20833333334; B = 378 }
    [PSCustomObject]@{ Row = 22; A = 45890.21875; B = 359 }
    [PSCustomObject]@{ Row = 23; A = 45890.22916666666; B = 345 }
    [PSCustomObject]@{ Row = 24; A = 45890.23958333334; B = 346 }
    [PSCustomObject]@{ Row = 25; A = 45890.25; B = 335 }
    [PSCustomObject]@{ Row = 26; A = 45890.26041666666; B = 331 }
    [PSCustomObject]@{ Row = 27; A = 45890.27083333334; B = 314 }
    [PSCustomObject]@{ Row = 28; A = 45890.28125; B = 292 }
    [PSCustomObject]@{ Row = 29; A = 45890.29166666666; B = 260 }
    [PSCustomObject]@{ Row = 30; A = 45890.30208333334; B = 223 }
    [PSCustomObject]@{ Row = 31; A = 45890.3125; B = 177 }
    [PSCustomObject]@{ Row = 32; A = 45890.32291666666; B = 150 }
    [PSCustomObject]@{ Row = 33; A = 45890.33333333334; B = 128 }
    [PSCustomObject]@{ Row = 34; A = 45890.34375; B = 107 }
    [PSCustomObject]@{ Row = 35; A = 45890.35416666666; B = 108 }
    [PSCustomObject]@{ Row = 36; A = 45890.36458333334; B = 111 }
    [PSCustomObject]@{ Row = 37; A = 45890.375; B = 107 }
    [PSCustomObject]@{ Row = 38; A = 45890.38541666666; B = 107 }
    [PSCustomObject]@{ Row = 39; A = 45890.39583333334; B = 113 }
    [PSCustomObject]@{ Row = 40; A = 45890.40625; B = 96 }
    [PSCustomObject]@{ Row = 41; A = 45890.41666666666; B = 0 }
    [PSCustomObject]@{ Row = 42; A = 45890.42708333334; B = 0 }
    [PSCustomObject]@{ Row = 43; A = 45890.4375; B = 0 }
    [PSCustomObject]@{ Row = 44; A = 45890.44791666666; B = 0 }
    [PSCustomObject]@{ Row = 45; A = 45890.45833333334; B = 0 }
    [PSCustomObject]@{ Row = 46; A = 45890.46875; B = 0 }
    [PSCustomObject]@{ Row = 47; A = 45890.47916666666; B = 0 }
    [PSCustomObject]@{ Row = 48; A = 45890.48958333334; B = 0 }
    [PSCustomObject]@{ Row = 49; A = 45890.5; B = 0 }
    [PSCustomObject]@{ Row = 50; A = 45890.51041666666; B = 0 }
    [PSCustomObject]@{ Row = 51; A = 45890.52083333334; B = 0 }
    [PSCustomObject]@{ Row = 52; A = 45890.53125; B = 0 }
    [PSCustomObject]@{ Row = 53; A = 45890.54166666666; B = 0 }
    [PSCustomObject]@{ Row = 54; A = 45890.55208333334; B = 0 }
    [PSCustomObject]@{ Row = 55; A = 45890.5625; B = 0 }
    [PSCustomObject]@{ Row = 56; A = 45890.57291666666; B = 0 }
    [PSCustomObject]@{ Row = 57; A = 45890.58333333334; B = 0 }
    [PSCustomObject]@{ Row = 58; A = 45890.59375; B = 0 }
    [PSCustomObject]@{ Row = 59; A = 45890.60416666666; B = 0 }
    [PSCustomObject]@{ Row = 60; A = 45890.61458333334; B = 0 }
    [PSCustomObject]@{ Row = 61; A = 45890.625; B = 0 }
    [PSCustomObject]@{ Row = 62; A = 45890.63541666666; B = 0 }
    [PSCustomObject]@{ Row = 63; A = 45890.64583333334; B = 0 }
    [PSCustomObject]@{ Row = 64; A = 45890.65625; B = 0 }
    [PSCustomObject]@{ Row = 65; A = 45890.66666666666; B = 0 }
    [PSCustomObject]@{ Row = 66; A = 45890.67708333334; B = 0 }
    [PSCustomObject]@{ Row = 67; A = 45890.6875; B = 0 }
    [PSCustomObject]@{ Row = 68; A = 45890.69791666666; B = 0 }
    [PSCustomObject]@{ Row = 69; A = 45890.70833333334; B = 0 }
    [PSCustomObject]@{ Row = 70; A = 45890.71875; B = 0 }
    [PSCustomObject]@{ Row = 71; A = 45890.72916666666; B = 0 }
    [PSCustomObject]@{ Row = 72; A = 45890.73958333334; B = 0 }
    [PSCustomObject]@{ Row = 73; A = 45890.75; B = 0 }
    [PSCustomObject]@{ Row = 74; A = 45890.76041666666; B = 0 }
    [PSCustomObject]@{ Row = 75; A = 45890.77083333334; B = 0 }
    [PSCustomObject]@{ Row = 76; A = 45890.78125; B = 0 }
    [PSCustomObject]@{ Row = 77; A = 45890.79166666666; B = 0 }
    [PSCustomObject]@{ Row = 78; A = 45890.80208333334; B = 0 }
    [PSCustomObject]@{ Row = 79; A = 45890.8125; B = 0 }
    [PSCustomObject]@{ Row = 80; A = 45890.82291666666; B = 0 }
    [PSCustomObject]@{ Row = 81; A = 45890.83333333334; B = 0 }
    [PSCustomObject]@{ Row = 82; A = 45890.84375; B = 0 }
    [PSCustomObject]@{ Row = 83; A = 45890.85416666666; B = 0 }
    [PSCustomObject]@{ Row = 84; A = 45890.86458333334; B = 0 }
    [PSCustomObject]@{ Row = 85; A = 45890.875; B = 0 }
    [PSCustomObject]@{ Row = 86; A = 45890.88541666666; B = 0 }
    [PSCustomObject]@{ Row = 87; A = 45890.89583333334; B = 0 }
    [PSCustomObject]@{ Row = 88; A = 45890.90625; B = 0 }
    [PSCustomObject]@{ Row = 89; A = 45890.91666666666; B = 0 }
    [PSCustomObject]@{ Row = 90; A = 45890.92708333334; B = 0 }
    [PSCustomObject]@{ Row = 91; A = 45890.9375; B = 0 }
    [PSCustomObject]@{ Row = 92; A = 45890.94791666666; B = 0 }
    [PSCustomObject]@{ Row = 93; A = 45890.95833333334; B = 0 }
    [PSCustomObject]@{ Row = 94; A = 45890.96875; B = 0 }
    [PSCustomObject]@{ Row = 95; A = 45890.97916666666; B = 0 }
    [PSCustomObject]@{ Row = 96; A = 45890.98958333334; B = 0 }
    [PSCustomObject]@{ Row = 97; A = 45891.0; B = 0 }
)

foreach ($item in $rowsData) {
    $ws.Cells.Item($item.Row, 1).Value2 = $item.A
    $ws.Cells.Item($item.Row, 2).Value2 = $item.B
}
